$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 125, shifting existing rows 125-245 down to 126-246.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new weekly record. Columns
# A,B,C,E,F,G,H,I,N,O,Q,R repeat the same values as the record now sitting at row
# 126 (the market/category/quality/unit metadata is identical week to week);
# D,J,K,L,M,P carry the new date/volume/price data per the diff.
$ws.Range("A125").Value = 4
$ws.Range("B125").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C125").Value = "Los Lagos"
$ws.Range("D125").Value = 44658
$ws.Range("E125").Value = 10
$ws.Range("F125").Value = 100112003
$ws.Range("G125").Value = "Ajo"
$ws.Range("H125").Value = "Chino"
$ws.Range("I125").Value = "Primera"
$ws.Range("J125").Value = 120
$ws.Range("K125").Value = 19000
$ws.Range("L125").Value = 20000
$ws.Range("M125").Value = 19500
$ws.Range("N125").Value = "$/caja 10 kilos"
$ws.Range("O125").Value = "China"
$ws.Range("P125").Value = 1950
$ws.Range("Q125").Value = 10
$ws.Range("R125").Value = "Hortaliza"
